{"js": "// Append a new \"Done\" bullet (same TextBody / ilvl=1 / numId=2 list item)\n// right after the last paragraph of the document body\n// (\"Done, please see \"testing/Test Plan.docx\"\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// insertParagraph inherits the paragraph formatting (style, numbering,\n// spacing, alignment, \u2026) of lastParagraph, matching the target diff which\n// reuses the same pPr as the preceding list item.\nlastParagraph.insertParagraph(\"Done\", \"After\");\n\nawait context.sync();\n", "ps1": "# Append a new \"Done\" bullet (same TextBody / ilvl=1 / numId=2 list item)\n# right after the last paragraph of the document\n# (\"Done, please see \"testing/Test Plan.docx\"\").\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n\n# InsertParagraphAfter clones the formatting (style, numbering, spacing,\n# alignment, \u2026) of the preceding paragraph, matching the target diff which\n# reuses the same pPr as the existing list item.\n$lastParagraph.Range.InsertParagraphAfter()\n\n$d.Paragraphs.Last.Range.Text = \"Done\"\n"}
